# Adjuste de titulo y grafico: agrega las nuevas filas de reservas
# capturadas (registros 56-65) a la hoja de datos "Sheet1".
#
# Las columnas C (Fecha Viaje), E (Cedula) y H (Telefono) contienen texto
# que "parece" una fecha o un numero (p.ej. "2025-03-21", "3333"). Si se
# asignan directamente via .Value, Excel las auto-convierte a fecha/numero.
# Para conservarlas como texto literal (igual que el resto de columnas de
# texto) forzamos el formato de celda a texto ("@") antes de escribir el
# valor y luego limpiamos el formato para que la celda quede con el estilo
# por defecto, igual que las filas existentes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$rows = @(
    @{ Row=56; A="Brasil";  B=1;  C="2025-03-21"; D="2025-03-21 13:33"; E="3333";     F="gg";     G="f";        H="3333";   I="f@gmail.cpm" },
    @{ Row=57; A="Dubai";   B=1;  C="2025-03-21"; D="2025-03-21 13:35"; E="555554";   F="lulu";   G="perez";    H="44";     I="p@gmail.com" },
    @{ Row=58; A="Francia"; B=1;  C="2025-03-21"; D="2025-03-21 13:41"; E="123";      F="lula";   G="li";       H="321";    I="lu@gmail.com" },
    @{ Row=59; A="España";  B=5;  C="2025-03-31"; D="2025-03-21 13:45"; E="122";      F="juan";   G="jun";      H="5656";   I="jj@gmai.com" },
    @{ Row=60; A="España";  B=10; C="2025-03-31"; D="2025-03-21 13:51"; E="100";      F="lina";   G="marin";    H="45545";  I="l@gmail.com" },
    @{ Row=61; A="España";  B=1;  C="2025-03-21"; D="2025-03-21 13:58"; E="444";      F="f";      G="ff";       H="444";    I="f@gmail.com" },
    @{ Row=62; A="España";  B=1;  C="2025-03-21"; D="2025-03-21 14:34"; E="23423423"; F="sfdsfs"; G="dsfsdfs";  H="234324"; I="wfwe@gmail.com" },
    @{ Row=63; A="España";  B=1;  C="2025-03-21"; D="2025-03-21 14:38"; E="444";      F="s";      G="s";        H="444";    I="s@gmal.com" },
    @{ Row=64; A="Dubai";   B=1;  C="2025-03-22"; D="2025-03-21 14:44"; E="3233333";  F="luis";   G="lopez";    H="333";    I="luis@gmai.com" },
    @{ Row=65; A="Dubai";   B=1;  C="2025-03-21"; D="2025-03-21 14:54"; E="32333444"; F="tulio";  G="lopez";    H="5555555";I="lopez@gmail.com" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    Set-TextValue "C$n" $r.C
    $ws.Range("D$n").Value = $r.D
    Set-TextValue "E$n" $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    Set-TextValue "H$n" $r.H
    $ws.Range("I$n").Value = $r.I
}
